# Generate Report for Handoff
# - Status changes from "Ready for handoff" to "Handoff transform failed"
# - The per-language "Latest Handoff File" hyperlink/value is cleared (handoff transform failed,
#   so no handoff file was produced) and the "Latest Handoff Datetime" / "Latest Handback DateTime"
#   fall back to the zero-date, with "Handoff Reason" becoming "Ignored".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

function Update-LanguageSheet($ws) {
    # Status: "Ready for handoff" -> "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # Remove the "Latest Handoff File" hyperlink + its value entirely (no handoff file produced).
    $target = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$C$2') {
            $target = $h
        }
    }
    if ($target -ne $null) {
        $target.Delete()
    }
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime falls back to the zero date.
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Latest Handback DateTime falls back to the zero date.
    $ws.Range("G2").Value = "0001-01-01 00:00:00"

    # Handoff Reason becomes "Ignored" instead of "Include".
    $ws.Range("H2").Value = "Ignored"
}

Update-LanguageSheet $ws2
Update-LanguageSheet $ws3

# Keep the Overview sheet's Status column in sync with the shared text.
$ws1.Range("B2").Value = "Handoff transform failed"
$ws1.Range("C2").Value = "Handoff transform failed"
